$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D as plain text so numeric-looking values (e.g. "24.00", "1.00")
# are not silently coerced into Number cells and lose their formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "58.128.43"
$ws.Range("E2").Value = "  +2.04%  "

# Row 3
$ws.Range("D3").Value = "2.359.70"
$ws.Range("E3").Value = "  +1.92%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "542.08"
$ws.Range("E5").Value = "  +2.33%  "

# Row 6
$ws.Range("D6").Value = "136.36"
$ws.Range("E6").Value = "  +2.93%  "

# Row 7
$ws.Range("E7").Value = "  +0.58%  "

# Row 8
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  +5.40%  "

# Row 9
$ws.Range("E9").Value = "  +1.79%  "

# Row 10
$ws.Range("D10").Value = "5.57"
$ws.Range("E10").Value = "  +3.92%  "

# Row 11
$ws.Range("E11").Value = "  -0.71%  "

# Row 12
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$ws.Range("D13").Value = "24.00"
$ws.Range("E13").Value = "  +2.72%  "

# Row 14
$ws.Range("D14").Value = "2.780.75"
$ws.Range("E14").Value = "  +1.70%  "

# Row 15
$ws.Range("D15").Value = "58.099.49"
$ws.Range("E15").Value = "  +1.98%  "

# Row 16
$ws.Range("D16").Value = "0.0000135"
$ws.Range("E16").Value = "  +1.57%  "

# Row 17
$ws.Range("D17").Value = "2.362.69"
$ws.Range("E17").Value = "  +1.56%  "

# Row 18
$ws.Range("D18").Value = "10.75"
$ws.Range("E18").Value = "  +3.16%  "

# Row 19
$ws.Range("D19").Value = "333.64"
$ws.Range("E19").Value = "  -0.59%  "

# Row 20
$ws.Range("E20").Value = "  +2.62%  "

# Row 21
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  -0.15%  "

# Row 22
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "62.80"
$ws.Range("E23").Value = "  +1.42%  "

# Row 24
$ws.Range("D24").Value = "0.167"
$ws.Range("E24").Value = "  +0.25%  "

# Row 25
$ws.Range("D25").Value = "8.53"
$ws.Range("E25").Value = "  -1.99%  "

# Row 26
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27
$ws.Range("E27").Value = "  +2.57%  "

# Row 28
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "1.75"
$ws.Range("E28").Value = "  +2.30%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "172.16"
$ws.Range("E29").Value = "  -0.56%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0742"
$ws.Range("E30").Value = "  +2.47%  "

# Row 31
$ws.Range("E31").Value = "  +1.03%  "

# Row 32
$ws.Range("E32").Value = "  +12.71%  "

# Row 33
$ws.Range("D33").Value = "18.53"
$ws.Range("E33").Value = "  +0.33%  "

# Row 34
$ws.Range("E34").Value = "  +0.05%  "

# Row 35
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +7.13%  "

# Row 36
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.75%  "

# Row 37
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  +3.98%  "

# Row 39
$ws.Range("D39").Value = "39.35"
$ws.Range("E39").Value = "  +0.52%  "

# Row 40
$ws.Range("D40").Value = "145.98"
$ws.Range("E40").Value = "  -1.78%  "

# Row 41
$ws.Range("D41").Value = "294.16"
$ws.Range("E41").Value = "  +3.72%  "

# Row 42
$ws.Range("D42").Value = "0.379"
$ws.Range("E42").Value = "  +1.16%  "

# Row 43
$ws.Range("E43").Value = "  +1.87%  "

# Row 44
$ws.Range("E44").Value = "  +1.75%  "

# Row 45
$ws.Range("D45").Value = "19.26"
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("E46").Value = "  +0.79%  "

# Row 47
$ws.Range("D47").Value = "0.565"
$ws.Range("E47").Value = "  +1.18%  "

# Row 48
$ws.Range("D48").Value = "0.0223"
$ws.Range("E48").Value = "  +2.88%  "

# Row 49
$ws.Range("D49").Value = "0.385"
$ws.Range("E49").Value = "  +1.01%  "

# Row 50
$ws.Range("D50").Value = "17.51"
$ws.Range("E50").Value = "  +0.47%  "

# Row 51
$ws.Range("D51").Value = "11.07"
$ws.Range("E51").Value = "  +0.50%  "
